# Auto-generated edit script: updates cached price/profit values
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match the refreshed Moogle market data snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 732.5833
$ws.Range("J88").Value = 817.4286
$ws.Range("L88").Value = 817.4286
$ws.Range("N88").Value = -1629.4286
$ws.Range("H91").Value = 732.5833
$ws.Range("J91").Value = 817.4286
$ws.Range("L91").Value = 817.4286
$ws.Range("N91").Value = -3625.4286
$ws.Range("H103").Value = 979.6667
$ws.Range("I103").Value = 1075.6
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 3226.8
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -2640.8
$ws.Range("N103").Value = -2672
$ws.Range("H106").Value = 40003132
$ws.Range("I106").Value = 48891420
$ws.Range("K106").Value = 48891420
$ws.Range("M106").Value = -48890789
$ws.Range("H132").Value = 2539.3064
$ws.Range("I132").Value = 1852.362
$ws.Range("K132").Value = 5557.086
$ws.Range("M132").Value = -3027.086
$ws.Range("H138").Value = 2293.5117
$ws.Range("I138").Value = 2168.5789
$ws.Range("J138").Value = 2392.4167
$ws.Range("K138").Value = 6505.736699999999
$ws.Range("L138").Value = 7177.250100000001
$ws.Range("M138").Value = -1365.736699999999
$ws.Range("N138").Value = -17457.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9056.645500000001
$ws.Range("I32").Value = 5772.864
$ws.Range("K32").Value = 5772.864
$ws.Range("M32").Value = -5485.864
$ws.Range("H45").Value = 4243.25
$ws.Range("I45").Value = 3831.8572
$ws.Range("K45").Value = 3831.8572
$ws.Range("M45").Value = -3454.8572
$ws.Range("H55").Value = 34179
$ws.Range("J55").Value = 46268.5
$ws.Range("L55").Value = 46268.5
$ws.Range("N55").Value = -46898.5
$ws.Range("H132").Value = 9171.333000000001
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H139").Value = 99628.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1829.4242
$ws.Range("I20").Value = 1378.8
$ws.Range("J20").Value = 2204.9443
$ws.Range("K20").Value = 1378.8
$ws.Range("L20").Value = 2204.9443
$ws.Range("M20").Value = -1131.8
$ws.Range("N20").Value = -2698.9443
$ws.Range("H94").Value = 6857.3335
$ws.Range("I94").Value = 7153.3335
$ws.Range("K94").Value = 7153.3335
$ws.Range("M94").Value = -6702.3335
$ws.Range("H107").Value = 1545.2941
$ws.Range("I107").Value = 1030.8334
$ws.Range("K107").Value = 1030.8334
$ws.Range("M107").Value = 889.1666
$ws.Range("H134").Value = 6188.7144
$ws.Range("I134").Value = 4664.2
$ws.Range("K134").Value = 13992.6
$ws.Range("M134").Value = -11457.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6742.803
$ws.Range("I31").Value = 4045.5881
$ws.Range("J31").Value = 10139.296
$ws.Range("K31").Value = 4045.5881
$ws.Range("L31").Value = 10139.296
$ws.Range("M31").Value = -3750.5881
$ws.Range("N31").Value = -10729.296
$ws.Range("H34").Value = 6742.803
$ws.Range("I34").Value = 4045.5881
$ws.Range("J34").Value = 10139.296
$ws.Range("K34").Value = 4045.5881
$ws.Range("L34").Value = 10139.296
$ws.Range("M34").Value = -3843.5881
$ws.Range("N34").Value = -10543.296
$ws.Range("H58").Value = 5836.875
$ws.Range("I58").Value = 4939
$ws.Range("K58").Value = 4939
$ws.Range("M58").Value = -4736
$ws.Range("H94").Value = 1214.909
$ws.Range("I94").Value = 800.25
$ws.Range("K94").Value = 800.25
$ws.Range("M94").Value = -349.25
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H134").Value = 3281.7297
$ws.Range("I134").Value = 2985.2188
$ws.Range("K134").Value = 8955.6564
$ws.Range("M134").Value = -6420.6564
$ws.Range("H136").Value = 5836.875
$ws.Range("I136").Value = 4939
$ws.Range("K136").Value = 14817
$ws.Range("M136").Value = -12267

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 641.2143
$ws.Range("I8").Value = 641.2143
$ws.Range("K8").Value = 1923.6429
$ws.Range("M8").Value = -1784.6429
$ws.Range("H11").Value = 375165.66
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H37").Value = 123333.336
$ws.Range("J37").Value = 123333.336
$ws.Range("L37").Value = 370000.008
$ws.Range("N37").Value = -370224.008
$ws.Range("H44").Value = 2089.4443
$ws.Range("J44").Value = 2770.5
$ws.Range("L44").Value = 8311.5
$ws.Range("N44").Value = -9107.5
$ws.Range("H46").Value = 4055
$ws.Range("I46").Value = 450
$ws.Range("J46").Value = 4570
$ws.Range("K46").Value = 1350
$ws.Range("L46").Value = 13710
$ws.Range("M46").Value = -1259
$ws.Range("N46").Value = -13892
$ws.Range("H120").Value = 16373.5
$ws.Range("I120").Value = 9253.223
$ws.Range("K120").Value = 27759.669
$ws.Range("M120").Value = -22921.669
$ws.Range("H123").Value = 7915.8335
$ws.Range("I123").Value = 831.6667
$ws.Range("J123").Value = 15000
$ws.Range("K123").Value = 2495.0001
$ws.Range("L123").Value = 45000
$ws.Range("M123").Value = -45.0001000000002
$ws.Range("N123").Value = -49900
$ws.Range("H131").Value = 3105.0476
$ws.Range("I131").Value = 1775.1818
$ws.Range("J131").Value = 4567.9
$ws.Range("K131").Value = 5325.5454
$ws.Range("L131").Value = 13703.7
$ws.Range("M131").Value = -285.5454
$ws.Range("N131").Value = -23783.7
$ws.Range("H140").Value = 1467.8182
$ws.Range("I140").Value = 1007.6667
$ws.Range("J140").Value = 2020
$ws.Range("K140").Value = 3023.0001
$ws.Range("L140").Value = 6060
$ws.Range("M140").Value = 2156.9999
$ws.Range("N140").Value = -16420

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 50000
$ws.Range("J53").Value = 50000
$ws.Range("L53").Value = 50000
$ws.Range("N53").Value = -51262
$ws.Range("H126").Value = 11118.667
$ws.Range("I126").Value = 11342.4
$ws.Range("K126").Value = 34027.2
$ws.Range("M126").Value = -31557.2
$ws.Range("H132").Value = 5060.875
$ws.Range("I132").Value = 3902.75
$ws.Range("K132").Value = 11708.25
$ws.Range("M132").Value = -9178.25
$ws.Range("H136").Value = 17931.592
$ws.Range("J136").Value = 17931.592
$ws.Range("L136").Value = 53794.776
$ws.Range("N136").Value = -58894.776

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 3494.6
$ws.Range("I46").Value = 2314
$ws.Range("J46").Value = 3789.75
$ws.Range("K46").Value = 2314
$ws.Range("L46").Value = 3789.75
$ws.Range("M46").Value = -2126
$ws.Range("N46").Value = -4165.75
$ws.Range("H55").Value = 1018.6667
$ws.Range("I55").Value = 126.42857
$ws.Range("J55").Value = 1586.4546
$ws.Range("K55").Value = 126.42857
$ws.Range("L55").Value = 1586.4546
$ws.Range("M55").Value = 46.57143000000001
$ws.Range("N55").Value = -1932.4546
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("H68").Value = 7653.5884
$ws.Range("I68").Value = 6109.1113
$ws.Range("K68").Value = 6109.1113
$ws.Range("M68").Value = -5360.1113
$ws.Range("H71").Value = 7653.5884
$ws.Range("I71").Value = 6109.1113
$ws.Range("K71").Value = 30545.5565
$ws.Range("M71").Value = -26801.5565
$ws.Range("H93").Value = 3129.9
$ws.Range("I93").Value = 1899.8572
$ws.Range("J93").Value = 6000
$ws.Range("K93").Value = 1899.8572
$ws.Range("L93").Value = 6000
$ws.Range("M93").Value = -651.8571999999999
$ws.Range("N93").Value = -8496
$ws.Range("H132").Value = 3292.375
$ws.Range("I132").Value = 1857.3928
$ws.Range("J132").Value = 13337.25
$ws.Range("K132").Value = 5572.178400000001
$ws.Range("L132").Value = 40011.75
$ws.Range("M132").Value = -3042.178400000001
$ws.Range("N132").Value = -45071.75
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 7806.6523
$ws.Range("I136").Value = 5514.8
$ws.Range("K136").Value = 16544.4
$ws.Range("M136").Value = -13994.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100000000
$ws.Range("I62").Value = 100000000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 100000000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -99999376
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 100000000
$ws.Range("I65").Value = 100000000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 500000000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -499996880
$ws.Range("N65").ClearContents()
$ws.Range("H75").Value = 82539
$ws.Range("I75").Value = 82539
$ws.Range("K75").Value = 82539
$ws.Range("M75").Value = -81603
$ws.Range("H78").Value = 82539
$ws.Range("I78").Value = 82539
$ws.Range("K78").Value = 247617
$ws.Range("M78").Value = -242937
